$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$texto = "Cidade registrou hoje alguns pontos de alagamento: previsão é de mais chuva. Repórter *ao vivo*. Choveu bastante na cidade. Trânsito ficou mais lento. Pista escorregadia. Alguns pontos ainda tem água. Muita atenção! Pode chover ainda mais no final da noite. Choveu forte em vários bairros. No Jóckey, algumas ruas ficaram alagadas. Rua Arnaldo Rosa Viana foi uma delas. Com previsão de mais chuva, IFF cancelou as aulas amanhã. fortes chuvas estavam previstas desde a noite de ontem. "

$rows = @(54, 55)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = "Record"
    $ws.Cells.Item($r, 2).Value = "RJ Record"
    $ws.Cells.Item($r, 3).Value = "Defesa Civil"
    $ws.Cells.Item($r, 4).Value = "2025-04-04T18:07"
    $ws.Cells.Item($r, 5).Value = "Neutro"
    $ws.Cells.Item($r, 6).Value = $texto
}
